$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine all product tuples from G2:G8 into a single Python-list-like string
$combined = "[('b.box led', '389.00'), ('Tomilho Serpilho', '10.00'), ('Oregano Bravo-Europeu', '10.00'), ('Manjericao Italiano', '10.00'), ('Coentro Portugues', '10.00'), ('Salsa Hortense', '10.00'), ('Alface Baby-Leaf', '10.00')]"

# Set the combined value into G2
$ws.Range("G2").Value = $combined

# Delete rows 3 through 8, since all products are now consolidated into row 2
$ws.Range("A3:G8").EntireRow.Delete()
